$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for the batch cohort generation columns
$ws.Range("G1").Value = "Dist"
$ws.Range("H1").Value = "CV"

# Populate the new columns for each parameter row (rows 2-9)
for ($row = 2; $row -le 9; $row++) {
    $ws.Cells.Item($row, 7).Value = "norm"
    $ws.Cells.Item($row, 8).Value = 0.05
}

# Keep the active selection in sync with the newly added column H
[void]$ws.Range("H12").Select()
